$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 141
    4  = 64
    6  = 135
    7  = 1296
    8  = 1549
    10 = 412
    11 = 182
    12 = 167
    13 = 164
    15 = 109
    16 = 273
    17 = 313
    18 = 327
    19 = 1751
    20 = 70
    23 = 680
    25 = 338
    26 = 4215
    28 = 280
    29 = 1104
    32 = 579
    34 = 278
    36 = 144
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
